# Break out stock.yaml completed
# - Fill in the previously-blank "backup" column (R) with the same value
#   as "detect_structure" (Q) for the rows where that backfill already
#   completed (rows 140, 160, 371, 399, 418, and the contiguous block
#   431-491).
# - Append 12 new trading days (rows 492-503) of OHLCV + derived data,
#   leaving their "backup" (R) column blank, same as the rows before the
#   backfill ran.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Backfill column R ("backup") from column Q ("detect_structure")
#    for the rows that were completed.
# ---------------------------------------------------------------------
$rowsToBackfill = @(140, 160, 371, 399, 418)
$rowsToBackfill += 431..491

foreach ($r in $rowsToBackfill) {
    $qVal = $ws.Cells.Item($r, 17).Value2   # column Q = 17
    $ws.Cells.Item($r, 18).Value2 = $qVal   # column R = 18
}

# ---------------------------------------------------------------------
# 2. Append the new rows of historical data (492-503).
#    Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close,
#    G Volume, H Year, I Month, J Day, K Hour, L Minute, M Second,
#    N Week, O isPivot, P two_line_structure, Q detect_structure.
#    Column R ("backup") is intentionally left blank for these rows.
# ---------------------------------------------------------------------
$newRows = @(
    @(492, 45630, 1607,                1610,                1589.650024414062, 1597.699951171875, 1597.699951171875, 553558,  2024, 12, 4,  0, 0, 0, 49, 0, 0, 0),
    @(493, 45631, 1610,                1620,                1561.25,           1577.900024414062, 1577.900024414062, 738747,  2024, 12, 5,  0, 0, 0, 49, 1, 0, 0),
    @(494, 45632, 1587.900024414062,   1587.900024414062,   1526.550048828125, 1547.849975585938, 1547.849975585938, 891217,  2024, 12, 6,  0, 0, 0, 49, 0, 0, 0),
    @(495, 45635, 1560,                1569,                1474.050048828125, 1480.099975585938, 1480.099975585938, 1228435, 2024, 12, 9,  0, 0, 0, 50, 0, 0, 0),
    @(496, 45636, 1489,                1495.449951171875,   1470,              1492.300048828125, 1492.300048828125, 665813,  2024, 12, 10, 0, 0, 0, 50, 0, 0, 0),
    @(497, 45637, 1500,                1511,                1477,              1486.599975585938, 1486.599975585938, 390322,  2024, 12, 11, 0, 0, 0, 50, 0, 0, 0),
    @(498, 45638, 1480,                1486.900024414062,   1455.150024414062, 1459.650024414062, 1459.650024414062, 419755,  2024, 12, 12, 0, 0, 0, 50, 0, 0, 0),
    @(499, 45639, 1460.800048828125,   1463.949951171875,   1438.099975585938, 1455.199951171875, 1455.199951171875, 541893,  2024, 12, 13, 0, 0, 0, 50, 0, 0, 0),
    @(500, 45642, 1455.300048828125,   1484.699951171875,   1455.300048828125, 1474.349975585938, 1474.349975585938, 476715,  2024, 12, 16, 0, 0, 0, 51, 0, 0, 0),
    @(501, 45643, 1481.449951171875,   1539,                1469.050048828125, 1494.599975585938, 1494.599975585938, 2061454, 2024, 12, 17, 0, 0, 0, 51, 0, 0, 0),
    @(502, 45644, 1485.25,             1489.5,              1404,              1408.449951171875, 1408.449951171875, 1866050, 2024, 12, 18, 0, 0, 0, 51, 0, 0, 0),
    @(503, 45645, 1395,                1425,                1385,              1421.25,            1421.25,           494881,  2024, 12, 19, 0, 0, 0, 51, 0, 0, 0)
)

# Number format used by the existing "Datetime" column cells.
$dateFormat = $ws.Cells.Item(491, 1).NumberFormat

foreach ($row in $newRows) {
    $r = $row[0]
    for ($col = 1; $col -le 17; $col++) {
        $ws.Cells.Item($r, $col).Value2 = $row[$col]
    }
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
}

Write-Output "edit complete"
